# Add PF/1.0.3 to meta-sheet
# Appends a new row (row 3) to Sheet1 with the version label "PF/1.0.3"
# and "X" markers in the sit2/uat2/prod columns, below the existing
# dev2/sit2/uat2/prod header (row 1) and PF/1.0.0 (row 2) rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 1).Value = "PF/1.0.3"
$ws.Cells.Item(3, 2).Value = "X"
$ws.Cells.Item(3, 3).Value = "X"
$ws.Cells.Item(3, 4).Value = "X"
